# Reproduce the authored changes to Placement.xlsx:
#  - Removal:                       move selection to G3
#  - FolioChildLocations:           move selection to E4 (and drop the stale topLeftCell)
#  - ChildLocationCorrespondences:  fill in J3/K3 and move selection to K3
#  - FolioChildLocationContactLogs: insert a new "NEW_CONTACT_LOG" column (F) with a
#                                    sample value, make this the active sheet, move
#                                    selection to H11
#  - ProviderSearch:                stops being the active sheet (handled implicitly by
#                                    activating FolioChildLocationContactLogs last)

$wb = $excel.ActiveWorkbook

# --- Removal -------------------------------------------------------------
$wsRemoval = $wb.Worksheets.Item("Removal")
$wsRemoval.Range("G3").Select()

# --- FolioChildLocations ---------------------------------------------------
$wsFolioChildLocations = $wb.Worksheets.Item("FolioChildLocations")
$wsFolioChildLocations.Range("E4").Select()

# --- ChildLocationCorrespondences ------------------------------------------
$wsCorrespondences = $wb.Worksheets.Item("ChildLocationCorrespondences")
$wsCorrespondences.Range("J3").Value = "past"
$wsCorrespondences.Range("K3").Value = "Court Ordered Placement"
$wsCorrespondences.Range("K3").Select()

# --- FolioChildLocationContactLogs -----------------------------------------
$wsContactLogs = $wb.Worksheets.Item("FolioChildLocationContactLogs")

# Insert a brand new column in front of the old column F, shifting the rest
# of the table (old F:AQ) one column to the right (new G:AR).
$wsContactLogs.Columns("F:F").Insert()

# Populate the new column's header + sample row.
$wsContactLogs.Range("F1").Value = "NEW_CONTACT_LOG"
$wsContactLogs.Range("F3").Value = "Click"

# This becomes the active sheet/tab, with the selection on H11.
$wsContactLogs.Activate()
$wsContactLogs.Range("H11").Select()
